# Scheduled-runner price/profit refresh: rewrites the market-price-derived
# columns (H:N = currentAveragePrice.., LevePriceNQ/HQ, LeveProfitNQ/HQ)
# for the leve rows whose underlying item prices moved since the last run.
# Column map: H=9999(current avg NQ) .. through N=LeveProfitHQ
#   8=H 9=I 10=J 11=K 12=L 13=M 14=N
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 9).Value = 3000
$ws.Cells.Item(113, 11).Value = 3000
$ws.Cells.Item(113, 13).Value = 254
$ws.Cells.Item(138, 8).Value = 1574.5143
$ws.Cells.Item(138, 9).Value = 1183.3658
$ws.Cells.Item(138, 10).Value = 2127.5173
$ws.Cells.Item(138, 11).Value = 3550.0974
$ws.Cells.Item(138, 12).Value = 6382.5519
$ws.Cells.Item(138, 13).Value = 1589.9026
$ws.Cells.Item(138, 14).Value = -16662.5519

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4301.75
$ws.Cells.Item(61, 9).Value = 5250
$ws.Cells.Item(61, 10).Value = 3985.6667
$ws.Cells.Item(61, 11).Value = 5250
$ws.Cells.Item(61, 12).Value = 3985.6667
$ws.Cells.Item(61, 13).Value = -5038
$ws.Cells.Item(61, 14).Value = -4409.6667
$ws.Cells.Item(82, 8).Value = 29999.25
$ws.Cells.Item(82, 10).Value = 29999.25
$ws.Cells.Item(82, 12).Value = 29999.25
$ws.Cells.Item(82, 14).Value = -30721.25
$ws.Cells.Item(85, 8).Value = 29999.25
$ws.Cells.Item(85, 10).Value = 29999.25
$ws.Cells.Item(85, 12).Value = 29999.25
$ws.Cells.Item(85, 14).Value = -32495.25
$ws.Cells.Item(122, 8).Value = 202962.4
$ws.Cells.Item(122, 9).Value = 335504
$ws.Cells.Item(122, 10).Value = 4150
$ws.Cells.Item(122, 11).Value = 1006512
$ws.Cells.Item(122, 12).Value = 12450
$ws.Cells.Item(122, 13).Value = -1004062
$ws.Cells.Item(122, 14).Value = -17350
$ws.Cells.Item(136, 8).Value = 4301.75
$ws.Cells.Item(136, 9).Value = 5250
$ws.Cells.Item(136, 10).Value = 3985.6667
$ws.Cells.Item(136, 11).Value = 15750
$ws.Cells.Item(136, 12).Value = 11957.0001
$ws.Cells.Item(136, 13).Value = -13200
$ws.Cells.Item(136, 14).Value = -17057.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 1878.8334
$ws.Cells.Item(99, 9).Value = 1548.6666
$ws.Cells.Item(99, 11).Value = 1548.6666
$ws.Cells.Item(99, 13).Value = -50.66660000000002
$ws.Cells.Item(126, 8).Value = 1878.8334
$ws.Cells.Item(126, 9).Value = 1548.6666
$ws.Cells.Item(126, 11).Value = 4645.9998
$ws.Cells.Item(126, 13).Value = -2175.9998
$ws.Cells.Item(132, 8).Value = 4066749.8
$ws.Cells.Item(132, 9).Value = 1279.6154
$ws.Cells.Item(132, 11).Value = 3838.8462
$ws.Cells.Item(132, 13).Value = -1308.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 8549.182
$ws.Cells.Item(62, 10).Value = 8549.182
$ws.Cells.Item(62, 12).Value = 25647.546
$ws.Cells.Item(62, 14).Value = -27019.546
$ws.Cells.Item(65, 8).Value = 8549.182
$ws.Cells.Item(65, 10).Value = 8549.182
$ws.Cells.Item(65, 12).Value = 76942.638
$ws.Cells.Item(65, 14).Value = -83806.638
$ws.Cells.Item(70, 8).Value = 2502.4
$ws.Cells.Item(70, 9).Value = 1506
$ws.Cells.Item(70, 10).Value = 3166.6667
$ws.Cells.Item(70, 11).Value = 4518
$ws.Cells.Item(70, 12).Value = 9500.000100000001
$ws.Cells.Item(70, 13).Value = -4203
$ws.Cells.Item(70, 14).Value = -10130.0001
$ws.Cells.Item(73, 8).Value = 2502.4
$ws.Cells.Item(73, 9).Value = 1506
$ws.Cells.Item(73, 10).Value = 3166.6667
$ws.Cells.Item(73, 11).Value = 4518
$ws.Cells.Item(73, 12).Value = 9500.000100000001
$ws.Cells.Item(73, 13).Value = -3426
$ws.Cells.Item(73, 14).Value = -11684.0001
$ws.Cells.Item(100, 8).Value = 2375
$ws.Cells.Item(100, 10).Value = 2375
$ws.Cells.Item(100, 12).Value = 7125
$ws.Cells.Item(100, 14).Value = -8747
$ws.Cells.Item(103, 8).Value = 951.2222
$ws.Cells.Item(103, 9).Value = 305.5
$ws.Cells.Item(103, 10).Value = 2242.6667
$ws.Cells.Item(103, 11).Value = 916.5
$ws.Cells.Item(103, 12).Value = 6728.000100000001
$ws.Cells.Item(103, 13).Value = -37.5
$ws.Cells.Item(103, 14).Value = -8486.000100000001
$ws.Cells.Item(106, 8).Value = 7946.143
$ws.Cells.Item(106, 9).Value = 2026
$ws.Cells.Item(106, 10).Value = 8932.833
$ws.Cells.Item(106, 11).Value = 6078
$ws.Cells.Item(106, 12).Value = 26798.499
$ws.Cells.Item(106, 13).Value = -5132
$ws.Cells.Item(106, 14).Value = -28690.499
$ws.Cells.Item(108, 8).Value = 885
$ws.Cells.Item(108, 9).Value = 885
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 11).Value = 2655
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 13).Value = 225
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(109, 8).Value = 3910
$ws.Cells.Item(109, 9).Value = 1175
$ws.Cells.Item(109, 10).Value = 5733.3335
$ws.Cells.Item(109, 11).Value = 3525
$ws.Cells.Item(109, 12).Value = 17200.0005
$ws.Cells.Item(109, 13).Value = -2485
$ws.Cells.Item(109, 14).Value = -19280.0005
$ws.Cells.Item(112, 8).Value = 3580
$ws.Cells.Item(112, 9).Value = 4750
$ws.Cells.Item(112, 10).Value = 2800
$ws.Cells.Item(112, 11).Value = 14250
$ws.Cells.Item(112, 12).Value = 8400
$ws.Cells.Item(112, 13).Value = -13142
$ws.Cells.Item(112, 14).Value = -10616
$ws.Cells.Item(115, 8).Value = 5725
$ws.Cells.Item(115, 9).Value = 3959.5
$ws.Cells.Item(115, 10).Value = 7490.5
$ws.Cells.Item(115, 11).Value = 11878.5
$ws.Cells.Item(115, 12).Value = 22471.5
$ws.Cells.Item(115, 13).Value = -10703.5
$ws.Cells.Item(115, 14).Value = -24821.5
$ws.Cells.Item(121, 8).Value = 445545.4
$ws.Cells.Item(121, 9).Value = 2000384.9
$ws.Cells.Item(121, 10).Value = 1305.5428
$ws.Cells.Item(121, 11).Value = 6001154.699999999
$ws.Cells.Item(121, 12).Value = 3916.6284
$ws.Cells.Item(121, 13).Value = -5999844.699999999
$ws.Cells.Item(121, 14).Value = -6536.6284
$ws.Cells.Item(122, 8).Value = 7947.2856
$ws.Cells.Item(122, 9).Value = 679.1429
$ws.Cells.Item(122, 11).Value = 6112.2861
$ws.Cells.Item(122, 13).Value = -3662.2861
$ws.Cells.Item(123, 8).Value = 5816.625
$ws.Cells.Item(123, 10).Value = 6361.857
$ws.Cells.Item(123, 12).Value = 19085.571
$ws.Cells.Item(123, 14).Value = -23985.571
$ws.Cells.Item(140, 8).Value = 1892.4375
$ws.Cells.Item(140, 9).Value = 1291.125
$ws.Cells.Item(140, 10).Value = 2493.75
$ws.Cells.Item(140, 11).Value = 3873.375
$ws.Cells.Item(140, 12).Value = 7481.25
$ws.Cells.Item(140, 13).Value = 1306.625
$ws.Cells.Item(140, 14).Value = -17841.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(33, 14).ClearContents()
$ws.Cells.Item(102, 8).Value = 1872
$ws.Cells.Item(102, 9).Value = 1834.3334
$ws.Cells.Item(102, 10).Value = 1928.5
$ws.Cells.Item(102, 11).Value = 1834.3334
$ws.Cells.Item(102, 12).Value = 1928.5
$ws.Cells.Item(102, 13).Value = -212.3334
$ws.Cells.Item(102, 14).Value = -5172.5
$ws.Cells.Item(122, 8).Value = 0
$ws.Cells.Item(122, 9).Value = 0
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).ClearContents()
$ws.Cells.Item(122, 14).ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2166.7058
$ws.Cells.Item(61, 9).Value = 1079.909
$ws.Cells.Item(61, 10).Value = 4159.1665
$ws.Cells.Item(61, 11).Value = 1079.909
$ws.Cells.Item(61, 12).Value = 4159.1665
$ws.Cells.Item(61, 13).Value = -877.9090000000001
$ws.Cells.Item(61, 14).Value = -4563.1665
$ws.Cells.Item(113, 8).Value = 2166.7058
$ws.Cells.Item(113, 9).Value = 1079.909
$ws.Cells.Item(113, 10).Value = 4159.1665
$ws.Cells.Item(113, 11).Value = 1079.909
$ws.Cells.Item(113, 12).Value = 4159.1665
$ws.Cells.Item(113, 13).Value = 1090.091
$ws.Cells.Item(113, 14).Value = -8499.1665
$ws.Cells.Item(122, 8).Value = 3084
$ws.Cells.Item(122, 9).Value = 3055.5557
$ws.Cells.Item(122, 10).Value = 3105.3333
$ws.Cells.Item(122, 11).Value = 9166.667099999999
$ws.Cells.Item(122, 12).Value = 9315.999899999999
$ws.Cells.Item(122, 13).Value = -6716.667099999999
$ws.Cells.Item(122, 14).Value = -14215.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2954.5454
$ws.Cells.Item(122, 9).Value = 2475
$ws.Cells.Item(122, 10).Value = 3530
$ws.Cells.Item(122, 11).Value = 7425
$ws.Cells.Item(122, 12).Value = 10590
$ws.Cells.Item(122, 13).Value = -4975
$ws.Cells.Item(122, 14).Value = -15490
$ws.Cells.Item(132, 8).Value = 5955276
$ws.Cells.Item(132, 9).Value = 4751.3
$ws.Cells.Item(132, 11).Value = 14253.9
$ws.Cells.Item(132, 13).Value = -11723.9
$ws.Cells.Item(136, 8).Value = 2143.4792
$ws.Cells.Item(136, 9).Value = 1846.4333
$ws.Cells.Item(136, 10).Value = 2638.5557
$ws.Cells.Item(136, 11).Value = 5539.2999
$ws.Cells.Item(136, 12).Value = 7915.6671
$ws.Cells.Item(136, 13).Value = -2989.2999
$ws.Cells.Item(136, 14).Value = -13015.6671
